$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the Metadata sheet: Version, Date and Contact values
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------------
# 2. The old "G1"/"G2" group values on the existing include sheets become
#    GUID identifiers; the plain "G1"/"G2" text moves onto two brand new
#    "Include from FSIII 3/4" sheets (copies of the existing ones).
# ---------------------------------------------------------------------------
$inc1 = $wb.Worksheets.Item("Include from FSIII")
$inc2 = $wb.Worksheets.Item("Include from FSIII 2")

$inc1.Range("C2").Value = "993d8f7b-fbed-4a78-90d9-6efbfa835114"
$inc2.Range("C2").Value = "ff47f955-3179-446f-b211-dc29de9456e3"

# ---------------------------------------------------------------------------
# 3. Add the two new sheets "Include from FSIII 3" and "Include from FSIII 4"
#    by copying "Include from FSIII 2" (to reuse its layout / styles) and
#    then setting the appropriate Value cell.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc2.Copy($null, $lastSheet)
$inc3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc3.Name = "Include from FSIII 3"
$inc3.Range("C2").Value = "G1"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc2.Copy($null, $lastSheet)
$inc4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc4.Name = "Include from FSIII 4"
$inc4.Range("C2").Value = "G2"

# Restore the originally active/selected sheet (Metadata) so the workbook
# view state is left unchanged, as in the source diff.
$meta.Select()
